$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM run added a third target cluster ("ECs") to the NATMI ligand-receptor
# table, so the original 2 sending x 2 target (FAPs, MuSCs) cross product (4 rows)
# becomes 2 sending x 3 target (ECs, FAPs, MuSCs) = 6 rows, with all the
# dependent specificity/weight metrics recomputed accordingly.
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T the various expression / detection-rate / specificity / edge-weight metrics.
$data = @(
    @("FAPs", "Ccl21b", "Ackr4", "ECs", 3, 1, 0.380435, 1.141305, 0.7997108917301441, 0.7997108917301442, 1, 0.3333333333333333, 0.07058066666666667, 0.211742, 0.03405243646146196, 0.03405243646146196, 0.02685135592333333, 0.24166220331, 0.02723210432817981, 0.02723210432817982),
    @("FAPs", "Ccl21b", "Ackr4", "FAPs", 3, 1, 0.380435, 1.141305, 0.7997108917301441, 0.7997108917301442, 3, 1, 1.890486333333333, 5.671459, 0.9120863940138783, 0.9120863940138783, 0.7192071682216667, 6.472864513995001, 0.7294054234917702, 0.7294054234917703),
    @("FAPs", "Ccl21b", "Ackr4", "MuSCs", 3, 1, 0.380435, 1.141305, 0.7997108917301441, 0.7997108917301442, 2, 0.6666666666666666, 0.1116383333333333, 0.334915, 0.05386116952465989, 0.05386116952465989, 0.04247112934166667, 0.382240164075, 0.04307336391019422, 0.04307336391019422),
    @("MuSCs", "Ccl21b", "Ackr4", "ECs", 2, 0.6666666666666666, 0.09528066666666667, 0.285842, 0.2002891082698559, 0.2002891082698559, 1, 0.3333333333333333, 0.07058066666666667, 0.211742, 0.03405243646146196, 0.03405243646146196, 0.006724972973777778, 0.060524756764, 0.006820332133282142, 0.006820332133282143),
    @("MuSCs", "Ccl21b", "Ackr4", "FAPs", 2, 0.6666666666666666, 0.09528066666666667, 0.285842, 0.2002891082698559, 0.2002891082698559, 3, 1, 1.890486333333333, 5.671459, 0.9120863940138783, 0.9120863940138783, 0.1801267981642222, 1.621141183478, 0.1826809705221081, 0.1826809705221081),
    @("MuSCs", "Ccl21b", "Ackr4", "MuSCs", 2, 0.6666666666666666, 0.09528066666666667, 0.285842, 0.2002891082698559, 0.2002891082698559, 2, 0.6666666666666666, 0.1116383333333333, 0.334915, 0.05386116952465989, 0.05386116952465989, 0.01063697482555556, 0.09573277343, 0.01078780561446567, 0.01078780561446567),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowVals[$j]
    }
}
